$d = $word.ActiveDocument

function Merge-ExactText {
    param(
        [string]$text
    )
    # Re-writes $text over itself via Find/Replace so that the runs
    # (and any w:proofErr spell-check markers) spanning it collapse
    # into a single run.
    $rng = $d.Content
    $null = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2)
}

function Resplit-Range {
    param(
        [int]$start,
        [int]$end
    )
    # Forces a run boundary at $start/$end (undoing any unwanted
    # coalescing with neighbouring runs that happen to share the same
    # formatting) by nudging the font size away and back again.
    $r = $d.Range($start, $end)
    $r.Font.Size = 99
    $r2 = $d.Range($start, $end)
    $r2.Font.Size = 12
}

# ---------------------------------------------------------------------
# 1) "Explorer, Google " + "Chrome" + ", Safari" -> one run
# ---------------------------------------------------------------------
Merge-ExactText "Explorer, Google Chrome, Safari"
$rng = $d.Content
$null = $rng.Find.Execute("Explorer, Google Chrome, Safari", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Resplit-Range $rng.Start $rng.End

# ---------------------------------------------------------------------
# 2) "El sistema " + "mobile" + " debe contar..." -> one run
# ---------------------------------------------------------------------
Merge-ExactText "El sistema mobile debe contar con textos que tengan un lenguaje amigable al usuario evitando el uso de términos técnicos"

# ---------------------------------------------------------------------
# 3) "Las contraseñas...forma " + "encriptada" -> one run
# ---------------------------------------------------------------------
Merge-ExactText "Las contraseñas se almacenaran en base de datos de forma encriptada"

# ---------------------------------------------------------------------
# 4) Fix typo "administrados" -> "administrador", split the sentence in
#    three runs and move the hidden _GoBack bookmark into the middle of
#    it (right after "de administrador").
# ---------------------------------------------------------------------
$rngFix = $d.Content
$null = $rngFix.Find.Execute( `
    "Únicamente los usuarios con rol de administrados podrán gestionar usuarios de la aplicación", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Únicamente los usuarios con rol de administrador podrán gestionar usuarios de la aplicación", 2)

$rngDe = $d.Content
$null = $rngDe.Find.Execute("de administrador", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$deStart = $rngDe.Start
$deEnd = $rngDe.End

# Split "...rol " | "de administrador" | " podrán..." into separate runs.
Resplit-Range $deStart $deEnd

# Relocate the _GoBack bookmark from wherever it currently sits to right
# after "de administrador".
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()
$bmRng = $d.Range($deEnd, $deEnd)
$null = $d.Bookmarks.Add("_GoBack", $bmRng)

# ---------------------------------------------------------------------
# 5) "El sistema  web debe contar con un sistema sencillo de " +
#    "deploy" + " en el servidor" -> one run
# ---------------------------------------------------------------------
Merge-ExactText "El sistema  web debe contar con un sistema sencillo de deploy en el servidor"

# ---------------------------------------------------------------------
# 6) "El proceso de desarrollo...y " + "deployarlos" + " en producción..." -> one run
# ---------------------------------------------------------------------
Merge-ExactText "El proceso de desarrollo debe permitir subir cambios al servidor web y deployarlos en producción de forma sencilla para minimizar tiempo de espera de resolución de fallos y de indisponibilidad del servicio"

# ---------------------------------------------------------------------
# 7) "Los sistemas tanto " + "mobile" + " como web...forma " -> one run,
#    keeping "homogénea" as its own trailing run.
# ---------------------------------------------------------------------
Merge-ExactText "Los sistemas tanto mobile como web deben contar con interfaces de usuarios definidas de forma homogénea"
$rngH = $d.Content
$null = $rngH.Find.Execute("homogénea", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Resplit-Range $rngH.Start $rngH.End
